$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-style rows 2:8 (cols A:B) with a thin box border, General number format ---
$dataRange = $ws.Range("A2:B8")
$dataRange.Style = "Normal"
$dataRange.Borders.LineStyle = 1

# --- Row 2 ---
$ws.Range("A2").Value = 13088334935
$ws.Range("B2").Value = 267794

# --- Row 3 ---
$ws.Range("A3").Value = 13098245418
$ws.Range("B3").Value = 313054

# --- Rows 4:8 now hold no data (cleared, but keep the new formatting) ---
$ws.Range("A4:B8").ClearContents()

# --- Drop the old row 9 entirely, shifting everything up ---
$ws.Rows.Item(9).Delete()

# --- Move the active selection to match the authored state ---
$null = $ws.Range("E27").Select()
